$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 229.8
$ws.Range("I92").Value = 199.76471
$ws.Range("K92").Value = 199.76471
$ws.Range("M92").Value = 1048.23529

# Row 96
$ws.Range("H96").Value = 333.73685
$ws.Range("I96").Value = 270.05884
$ws.Range("K96").Value = 810.17652
$ws.Range("M96").Value = 562.82348

# Row 99
$ws.Range("H99").Value = 948.125
$ws.Range("I99").Value = 411.2
$ws.Range("J99").Value = 1843
$ws.Range("K99").Value = 1233.6
$ws.Range("L99").Value = 5529
$ws.Range("M99").Value = 264.4000000000001
$ws.Range("N99").Value = -8525

# Row 107
$ws.Range("H107").Value = 966.6818
$ws.Range("I107").Value = 932.3125
$ws.Range("K107").Value = 932.3125
$ws.Range("M107").Value = 987.6875

# Row 137
$ws.Range("H137").Value = 2745.238
$ws.Range("I137").Value = 1070.8889
$ws.Range("J137").Value = 4001
$ws.Range("K137").Value = 3212.6667
$ws.Range("L137").Value = 12003
$ws.Range("M137").Value = -662.6666999999998
$ws.Range("N137").Value = -17103

# Row 138
$ws.Range("H138").Value = 2767.82
$ws.Range("I138").Value = 1145.683
$ws.Range("J138").Value = 3895.0679
$ws.Range("K138").Value = 3437.049
$ws.Range("L138").Value = 11685.2037
$ws.Range("M138").Value = 1702.951
$ws.Range("N138").Value = -21965.2037


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23066.123
$ws.Range("I32").Value = 8558.65
$ws.Range("K32").Value = 8558.65
$ws.Range("M32").Value = -8271.65

# Row 61
$ws.Range("H61").Value = 2387.9614
$ws.Range("I61").Value = 1999.409
$ws.Range("J61").Value = 4525
$ws.Range("K61").Value = 1999.409
$ws.Range("L61").Value = 4525
$ws.Range("M61").Value = -1787.409
$ws.Range("N61").Value = -4949

# Row 74
$ws.Range("H74").Value = 24358.844
$ws.Range("I74").Value = 1371.2354
$ws.Range("J74").Value = 95411.45
$ws.Range("K74").Value = 1371.2354
$ws.Range("L74").Value = 95411.45
$ws.Range("M74").Value = -497.2354
$ws.Range("N74").Value = -97159.45

# Row 77
$ws.Range("H77").Value = 24358.844
$ws.Range("I77").Value = 1371.2354
$ws.Range("J77").Value = 95411.45
$ws.Range("K77").Value = 6856.177
$ws.Range("L77").Value = 477057.25
$ws.Range("M77").Value = -2488.177
$ws.Range("N77").Value = -485793.25

# Row 132
$ws.Range("H132").Value = 2383.7222
$ws.Range("I132").Value = 1904.9131
$ws.Range("J132").Value = 3230.8462
$ws.Range("K132").Value = 5714.7393
$ws.Range("L132").Value = 9692.5386
$ws.Range("M132").Value = -3184.7393
$ws.Range("N132").Value = -14752.5386

# Row 136
$ws.Range("H136").Value = 2387.9614
$ws.Range("I136").Value = 1999.409
$ws.Range("J136").Value = 4525
$ws.Range("K136").Value = 5998.227000000001
$ws.Range("L136").Value = 13575
$ws.Range("M136").Value = -3448.227000000001
$ws.Range("N136").Value = -18675


$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3091.3
$ws.Range("I134").Value = 2328.4285
$ws.Range("J134").Value = 4871.3335
$ws.Range("K134").Value = 6985.2855
$ws.Range("L134").Value = 14614.0005
$ws.Range("M134").Value = -4450.2855
$ws.Range("N134").Value = -19684.0005


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4461.641
$ws.Range("I31").Value = 2161.7896
$ws.Range("J31").Value = 6646.5
$ws.Range("K31").Value = 2161.7896
$ws.Range("L31").Value = 6646.5
$ws.Range("M31").Value = -1866.7896
$ws.Range("N31").Value = -7236.5

# Row 34
$ws.Range("H34").Value = 4461.641
$ws.Range("I34").Value = 2161.7896
$ws.Range("J34").Value = 6646.5
$ws.Range("K34").Value = 2161.7896
$ws.Range("L34").Value = 6646.5
$ws.Range("M34").Value = -1959.7896
$ws.Range("N34").Value = -7050.5

# Row 58
$ws.Range("H58").Value = 1774.8
$ws.Range("I58").Value = 1564.4706
$ws.Range("J58").Value = 2966.6667
$ws.Range("K58").Value = 1564.4706
$ws.Range("L58").Value = 2966.6667
$ws.Range("M58").Value = -1361.4706
$ws.Range("N58").Value = -3372.6667

# Row 107
$ws.Range("H107").Value = 820.1667
$ws.Range("J107").Value = 648.1
$ws.Range("L107").Value = 648.1
$ws.Range("N107").Value = -4488.1

# Row 134
$ws.Range("H134").Value = 3425.568
$ws.Range("I134").Value = 3643.275
$ws.Range("J134").Value = 1248.5
$ws.Range("K134").Value = 10929.825
$ws.Range("L134").Value = 3745.5
$ws.Range("M134").Value = -8394.825000000001
$ws.Range("N134").Value = -8815.5

# Row 136
$ws.Range("H136").Value = 1774.8
$ws.Range("I136").Value = 1564.4706
$ws.Range("J136").Value = 2966.6667
$ws.Range("K136").Value = 4693.4118
$ws.Range("L136").Value = 8900.000100000001
$ws.Range("M136").Value = -2143.4118
$ws.Range("N136").Value = -14000.0001


$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 2254.85
$ws.Range("I109").Value = 706.0909
$ws.Range("J109").Value = 4147.778
$ws.Range("K109").Value = 2118.2727
$ws.Range("L109").Value = 12443.334
$ws.Range("M109").Value = -1078.2727
$ws.Range("N109").Value = -14523.334

# Row 119
$ws.Range("H119").Value = 3824.5881
$ws.Range("I119").Value = 2693.6924
$ws.Range("J119").Value = 7500
$ws.Range("K119").Value = 8081.0772
$ws.Range("L119").Value = 22500
$ws.Range("M119").Value = -3243.0772
$ws.Range("N119").Value = -32176

# Row 120
$ws.Range("H120").Value = 10857.857
$ws.Range("I120").Value = 2851.6667
$ws.Range("J120").Value = 16862.5
$ws.Range("K120").Value = 8555.000100000001
$ws.Range("L120").Value = 50587.5
$ws.Range("M120").Value = -3717.000100000001
$ws.Range("N120").Value = -60263.5

# Row 129
$ws.Range("H129").Value = 2003.75
$ws.Range("I129").Value = 3030
$ws.Range("J129").Value = 1857.1428
$ws.Range("K129").Value = 9090
$ws.Range("L129").Value = 5571.428400000001
$ws.Range("M129").Value = -4090
$ws.Range("N129").Value = -15571.4284

# Row 131
$ws.Range("H131").Value = 912.4299999999999
$ws.Range("J131").Value = 950.25
$ws.Range("L131").Value = 2850.75
$ws.Range("N131").Value = -12930.75

# Row 136
$ws.Range("H136").Value = 3833.422
$ws.Range("I136").Value = 1872.2222
$ws.Range("J136").Value = 4051.3333
$ws.Range("K136").Value = 5616.6666
$ws.Range("L136").Value = 12153.9999
$ws.Range("M136").Value = -516.6665999999996
$ws.Range("N136").Value = -22353.9999


$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 7666.467
$ws.Range("I132").Value = 12334
$ws.Range("J132").Value = 4554.778
$ws.Range("K132").Value = 37002
$ws.Range("L132").Value = 13664.334
$ws.Range("M132").Value = -34472
$ws.Range("N132").Value = -18724.334

# Row 136
$ws.Range("H136").Value = 3808.9546
$ws.Range("I136").Value = 1407.5405
$ws.Range("K136").Value = 4222.6215
$ws.Range("M136").Value = -1672.6215


$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 38107.7
$ws.Range("J82").Value = 38107.7
$ws.Range("L82").Value = 38107.7
$ws.Range("N82").Value = -38873.7

# Row 85
$ws.Range("H85").Value = 38107.7
$ws.Range("J85").Value = 38107.7
$ws.Range("L85").Value = 38107.7
$ws.Range("N85").Value = -40759.7

